$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 20 - "Exceptions"
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Add(20, 2)
$s20.Shapes.Item(1).TextFrame.TextRange.Text = "Exceptions"

$body20 = $s20.Shapes.Item(2).TextFrame.TextRange
$body20.Text = "try{stuff}catch(err){"
$body20.InsertAfter("console.log") | Out-Null
$body20.InsertAfter("(err);}finally{do anyway}") | Out-Null

# ---------------------------------------------------------------------------
# Slide 21 - "Debugging"
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Add(21, 2)
$s21.Shapes.Item(1).TextFrame.TextRange.Text = "Debugging"

$body21 = $s21.Shapes.Item(2).TextFrame.TextRange
$body21.Text = "Console.log"
$body21.InsertAfter("() ;)") | Out-Null
$body21.InsertAfter("`rdebugger; creates a ") | Out-Null
$body21.InsertAfter("Breackpoint") | Out-Null

# ---------------------------------------------------------------------------
# Slide 22 - "Your Turn"
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Add(22, 2)
$s22.Shapes.Item(1).TextFrame.TextRange.Text = "Your Turn"

$body22 = $s22.Shapes.Item(2).TextFrame.TextRange
$body22.Text = "Install Editor of your choice (VS "
$body22.InsertAfter("Code)") | Out-Null
$body22.InsertAfter("`rOpen ") | Out-Null
$body22.InsertAfter("Code Academy") | Out-Null
$body22.InsertAfter("`rCreate Account and start with exercises") | Out-Null
$codeAcademy = $body22.Characters(46, 12)
$codeAcademy.ActionSettings(1).Hyperlink.Address = "https://www.codecademy.com"

# ---------------------------------------------------------------------------
# Slide 23 - "Your turn, again"
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Add(23, 2)
$s23.Shapes.Item(1).TextFrame.TextRange.Text = "Your turn, again"

$body23 = $s23.Shapes.Item(2).TextFrame.TextRange
$body23.Text = "Create "
$body23.InsertAfter("SlideShow") | Out-Null
$body23.InsertAfter(" like ") | Out-Null
$body23.InsertAfter("here") | Out-Null
$body23.InsertAfter(" but with JavaScript") | Out-Null
$here = $body23.Characters(23, 4)
$here.ActionSettings(1).Hyperlink.Address = "https://revealjs.com"
